$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date headers for columns AW:AZ (continuing the daily series after 14-ago)
$ws.Range("AW1").Value = "16-ago"
$ws.Range("AX1").Value = "17-ago"
$ws.Range("AY1").Value = "18-ago"
$ws.Range("AZ1").Value = "22-ago"

# Copy the header formatting (text number format) from the previous header cell
$ws.Range("AV1").Copy()
$ws.Range("AW1:AZ1").PasteSpecial(-4122)

# New data values for rows 2:11, columns AW:AZ
$rows = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11)
$aw = @(16, 14, 10, 12, 13, 10, 10, 20, 14, 23)
$ax = @(15, 13, 10, 13, 11, 11, 10, 19, 13, 22)
$ay = @(16, 14, 11, 13, 10, 11, 11, 18, 11, 20)
$az = @(16, 13, 11, 11, 13, 11, 14, 16, 7, 16)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $ws.Cells.Item($r, 49).Value = $aw[$i]
    $ws.Cells.Item($r, 50).Value = $ax[$i]
    $ws.Cells.Item($r, 51).Value = $ay[$i]
    $ws.Cells.Item($r, 52).Value = $az[$i]
}

# Copy the data cell formatting (integer, centered) from the previous data column
$ws.Range("AV2").Copy()
$ws.Range("AW2:AZ11").PasteSpecial(-4122)

# Extend the hidden/zero-width column formatting that previously covered J:AL (10:38)
# so that it also covers the newly filled-in AM:AT (39:46) range
$ws.Range("AM1:AT1").EntireColumn.ColumnWidth = -0.8333333333333334
$ws.Range("AM1:AT1").EntireColumn.Hidden = $true

# Update the active selection as recorded after the edit
$ws.Range("BC8").Select()
